$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on D:E columns (rows 2-51) before writing, to prevent Excel
# from auto-converting numeric-looking strings into real numbers, then restore
# the default "Normal" style so the cell style index matches the original file.
$numRange = $ws.Range("D2:E51")
$numRange.NumberFormat = "@"

$ws.Range("D2").Value = '25.707.63'
$ws.Range("E2").Value = '  -3.18%  '
$ws.Range("D3").Value = '1.761.78'
$ws.Range("E3").Value = '  -4.34%  '
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.74%  '
$ws.Range("D5").Value = '233.78'
$ws.Range("E5").Value = '  -10.01%  '
$ws.Range("D6").Value = '1.006'
$ws.Range("E6").Value = '  +0.53%  '
$ws.Range("D7").Value = '0.4875'
$ws.Range("E7").Value = '  -7.01%  '
$ws.Range("D8").Value = '42.19'
$ws.Range("E8").Value = '  -6.80%  '
$ws.Range("D9").Value = '0.2406'
$ws.Range("E9").Value = '  -24.38%  '
$ws.Range("D10").Value = '0.05970'
$ws.Range("E10").Value = '  -12.16%  '
$ws.Range("D11").Value = '1.776.21'
$ws.Range("E11").Value = '  -1.66%  '
$ws.Range("D12").Value = '0.06603'
$ws.Range("E12").Value = '  -15.00%  '
$ws.Range("D13").Value = '13.66'
$ws.Range("E13").Value = '  -27.15%  '
$ws.Range("D14").Value = '0.5837'
$ws.Range("E14").Value = '  -25.28%  '
$ws.Range("D15").Value = '76.17'
$ws.Range("E15").Value = '  -13.43%  '
$ws.Range("D16").Value = '4.266'
$ws.Range("E16").Value = '  -14.96%  '
$ws.Range("D17").Value = '1.010'
$ws.Range("E17").Value = '  +0.93%  '
$ws.Range("E18").Value = '  +0.37%  '
$ws.Range("D19").Value = '25.712.88'
$ws.Range("E19").Value = '  -3.24%  '
$ws.Range("D20").Value = '10.77'
$ws.Range("E20").Value = '  -22.36%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '1.994.25'
$ws.Range("E21").Value = '  -3.27%  '
$ws.Range("B22").Value = 'ShibaInu'
$ws.Range("C22").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D22").Value = '0.000006173'
$ws.Range("E22").Value = '  -22.30%  '
$ws.Range("D23").Value = '3.794'
$ws.Range("E23").Value = '  -17.71%  '
$ws.Range("D24").Value = '5.058'
$ws.Range("E24").Value = '  -15.34%  '
$ws.Range("D25").Value = '7.897'
$ws.Range("E25").Value = '  -15.51%  '
$ws.Range("D26").Value = '132.21'
$ws.Range("E26").Value = '  -7.38%  '
$ws.Range("D27").Value = '1.842'
$ws.Range("E27").Value = '  -17.19%  '
$ws.Range("D28").Value = '1.408'
$ws.Range("E28").Value = '  -15.97%  '
$ws.Range("D29").Value = '14.11'
$ws.Range("E29").Value = '  -16.46%  '
$ws.Range("D30").Value = '98.61'
$ws.Range("E30").Value = '  -11.89%  '
$ws.Range("D31").Value = '0.08199'
$ws.Range("E31").Value = '  -6.10%  '
$ws.Range("D32").Value = '3.558'
$ws.Range("E32").Value = '  -14.97%  '
$ws.Range("B33").Value = 'Frax'
$ws.Range("C33").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D33").Value = '1.004'
$ws.Range("E33").Value = '  +0.72%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '0.04221'
$ws.Range("E34").Value = '  -13.65%  '
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").Value = '3.110'
$ws.Range("E35").Value = '  -23.70%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = '2.618'
$ws.Range("E36").Value = '  -8.35%  '
$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").Value = '1.012'
$ws.Range("E37").Value = '  -10.67%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = '0.5995'
$ws.Range("E38").Value = '  -17.20%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.710'
$ws.Range("E39").Value = '  -12.37%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = '2.068'
$ws.Range("E40").Value = '  -7.30%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").Value = '1.005'
$ws.Range("E41").Value = '  +0.45%  '
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").Value = '100.78'
$ws.Range("E42").Value = '  -8.66%  '
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").Value = '0.01432'
$ws.Range("E43").Value = '  -17.93%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").Value = '0.7759'
$ws.Range("E44").Value = '  -13.62%  '
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").Value = '5.152'
$ws.Range("E45").Value = '  -13.01%  '
$ws.Range("B46").Value = 'TheSandbox'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D46").Value = '0.3741'
$ws.Range("E46").Value = '  -22.44%  '
$ws.Range("D47").Value = '0.05139'
$ws.Range("E47").Value = '  -11.97%  '
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = '6.021'
$ws.Range("E48").Value = '  -21.28%  '
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").Value = '51.72'
$ws.Range("E49").Value = '  -13.21%  '
$ws.Range("B50").Value = 'USDD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D50").Value = '1.004'
$ws.Range("E50").Value = '  -0.04%  '
$ws.Range("B51").Value = 'TrueUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range("D51").Value = '1.003'
$ws.Range("E51").Value = '  +0.33%  '

$numRange.Style = "Normal"
